$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Weekly Quantity": append new row 25
# ---------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("A25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeekly.Range("A25").Value2 = 45662.99999999999
$wsWeekly.Range("B25").Value2 = 1

# ---------------------------------------------------------------
# Sheet "Monthly Trend": append new row 10
# ---------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsMonthly.Range("A10").Value2 = 45688.99999999999
$wsMonthly.Range("B10").Value2 = 1

# ---------------------------------------------------------------
# Sheet "PO Forecast": new forecast model values
# ---------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Updated forecast quantities for existing rows (dates unchanged)
$wsForecast.Range("B2").Value2 = 16
$wsForecast.Range("B3").Value2 = 16
$wsForecast.Range("B4").Value2 = 16
$wsForecast.Range("B8").Value2 = 15
$wsForecast.Range("B14").Value2 = 12
$wsForecast.Range("B15").Value2 = 12
$wsForecast.Range("B17").Value2 = 11
$wsForecast.Range("B18").Value2 = 11
$wsForecast.Range("B19").Value2 = 11
$wsForecast.Range("B20").Value2 = 11
$wsForecast.Range("B21").Value2 = 10
$wsForecast.Range("B22").Value2 = 10
$wsForecast.Range("B23").Value2 = 9
$wsForecast.Range("B24").Value2 = 9

# Rows 25-32 shift to later dates with new quantities
$wsForecast.Range("A25").Value2 = 45662.99999999999
$wsForecast.Range("B25").Value2 = 8

$wsForecast.Range("A26").Value2 = 45669.99999999999
$wsForecast.Range("B26").Value2 = 8

$wsForecast.Range("A27").Value2 = 45676.99999999999
$wsForecast.Range("B27").Value2 = 8

$wsForecast.Range("A28").Value2 = 45683.99999999999
$wsForecast.Range("B28").Value2 = 8

$wsForecast.Range("A29").Value2 = 45690.99999999999
$wsForecast.Range("B29").Value2 = 7

$wsForecast.Range("A30").Value2 = 45697.99999999999
$wsForecast.Range("B30").Value2 = 7

$wsForecast.Range("A31").Value2 = 45704.99999999999
$wsForecast.Range("B31").Value2 = 7

$wsForecast.Range("A32").Value2 = 45711.99999999999
$wsForecast.Range("B32").Value2 = 7

# New row 33
$wsForecast.Range("A33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("A33").Value2 = 45718.99999999999
$wsForecast.Range("B33").Value2 = 6
